# fix(gui) step 1 and 2
#
# The sheet "CAJAS PARA ENRROLLADOR" keeps a running "last updated" date in
# A1 (stored as a date serial) plus two prices (step 1 / step 2) in D30 and
# D31. Bump the date by a day and update both prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Date in A1: 17-Jan-2024 (45308) -> 18-Jan-2024 (45309)
$ws.Range("A1").Value = 45309

# Step 1 price (CAJA p/ ENROLLADOR CHICA)
$ws.Range("D30").Value = 570

# Step 2 price (CAJA p/ ENROLLADOR GRANDE)
$ws.Range("D31").Value = 690
